$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2.9
$ws.Range("J2").Value = 3.25
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("W2").Value = 8
$ws.Range("AL2").Value = 23
